# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
# F4: 1468 -> 1469
# F6: 27   -> 28
# F9: 257  -> 258

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1469
    $ws.Range("F6").Value = 28
    $ws.Range("F9").Value = 258
}
